# Salvando trabalho local antes do merge
#
# Refreshes the two sample data rows with new placeholder/test values and
# reinserts two columns next to their related fields:
#   - prestador_inscricao_municipal (now right after prestador_razao_social)
#   - tomador_email (now right after tomador_razao_social)
# which shifts the other prestador_*/tomador_* address columns over by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage. A few sample values
# are purely numeric strings (CNPJ/CEP/document-number style codes) that
# must stay text instead of being coerced into numbers by Excel's normal
# auto-detection. NumberFormat is reset back to the default afterwards so
# no stray cell formatting is left behind.
function Set-TextValue($sheet, $addr, $val) {
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $val
    $sheet.Range($addr).Style = "Normal"
}

# --- Row 1 headers: columns shift because two columns were reinserted ----
$ws.Range("I1").Value = "prestador_inscricao_municipal"
$ws.Range("J1").Value = "prestador_logradouro"
$ws.Range("K1").Value = "prestador_bairro"
$ws.Range("L1").Value = "prestador_cep"
$ws.Range("M1").Value = "prestador_cidade"
$ws.Range("N1").Value = "prestador_uf"
$ws.Range("Q1").Value = "tomador_email"
$ws.Range("R1").Value = "tomador_logradouro"
$ws.Range("S1").Value = "tomador_bairro"
$ws.Range("T1").Value = "tomador_cep"
$ws.Range("U1").Value = "tomador_cidade"
$ws.Range("V1").Value = "tomador_uf"

# --- Row 2 -----------------------------------------------------------
$ws.Range("A2").Value = "6b4ccb05496145fda961038a16b1a2d1"
$ws.Range("B2").Value = "nfe.jpg"
$ws.Range("C2").Value = "2025-09-14 14:39:04"
Set-TextValue $ws "D2" "0001234567890123"
$ws.Range("E2").Value = "2023-03-15T12:34:56.789"
Set-TextValue $ws "F2" "1234567890123456"
$ws.Range("G2").Value = "12.345.678/00-00"
$ws.Range("H2").Value = "SERVICOS DE ELECTRONICA LTDA."
$ws.Range("I2").Value = "INSCRIÇÃO MUNICIPAL 123456789012345"
$ws.Range("J2").Value = "Rua Exemplo, 123 - Bairro Novo, Cidade"
$ws.Range("K2").Value = "Novos"
$ws.Range("L2").Value = "12345-678"
$ws.Range("M2").Value = "Cidade"
$ws.Range("N2").Value = "SP"
$ws.Range("O2").Value = "00.000.000-00"
$ws.Range("P2").Value = "EMPRESA DE COMÉRCIO ESTÁTICO"
$ws.Range("Q2").Value = "contato@exemplo.com"
$ws.Range("R2").Value = "Rua Exemplo, 456 - Bairro Novo, Cidade"
$ws.Range("S2").Value = "Novo"
$ws.Range("T2").Value = "12345-000"
$ws.Range("U2").Value = "Cidade"
$ws.Range("V2").Value = "SP"
$ws.Range("W2").Value = "TROCA DE SERVICO"
$ws.Range("X2").Value = "99.999.000-00"
$ws.Range("Y2").Value = "SERVIÇO DE TROCA DE ELETRÔNICOS"
$ws.Range("Z2").Value = 600
$ws.Range("AA2").Value = 500
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 90
$ws.Range("AD2").Value = 70

# --- Row 3 -----------------------------------------------------------
$ws.Range("A3").Value = "405a233a317d0adf17b0f4d02beef0e5"
$ws.Range("B3").Value = "WhatsApp Image 2025-08-20 at 20.50.35.jpeg"
$ws.Range("C3").Value = "2025-09-14 14:40:03"
Set-TextValue $ws "D3" "0000000000000000"
$ws.Range("E3").Value = "2019-09-30T00:00:00"
Set-TextValue $ws "F3" "1234567890"
$ws.Range("G3").Value = "12.345.678/00-00"
$ws.Range("H3").Value = "SERVICOS DE CONSULTORIA LTDA."
Set-TextValue $ws "I3" "1234567890"
$ws.Range("J3").Value = "Rua Exemplo, 123 - Bairro Novo, Cidade, Estado"
$ws.Range("K3").Value = "Novos"
Set-TextValue $ws "L3" "12345678"
$ws.Range("M3").Value = "Cidade"
$ws.Range("N3").Value = "SP"
$ws.Range("O3").Value = "00.000.000-00"
$ws.Range("P3").Value = "Nome da Empresa"
$ws.Range("Q3").Value = "contato@exemplo.com"
$ws.Range("R3").Value = "Rua Exemplo, 234 - Bairro Novo, Cidade, Estado"
$ws.Range("S3").Value = "Novo"
Set-TextValue $ws "T3" "12345678"
$ws.Range("U3").Value = "Cidade"
$ws.Range("V3").Value = "SP"
$ws.Range("W3").Value = "SERVICOS DE CONSULTORIA"
$ws.Range("X3").Value = "99.000.000-00"
$ws.Range("Y3").Value = "Consultoria em Marketing"
$ws.Range("Z3").Value = 1500
$ws.Range("AA3").Value = 1500
$ws.Range("AB3").Value = 0
$ws.Range("AC3").Value = 1575
$ws.Range("AD3").Value = 75
